$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")

# Header updates
$ws1.Cells.Item(2, 1).Value = "Última actualización: 11:20:08"
$ws1.Cells.Item(3, 1).Value = "Total filas: 200"

# Row data updates/additions (Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada)
$ws1Rows = @{
    42 = @('05:18:56', '06:58', '10_OLMOS', 100, 'LP1912')
    43 = @('06:57:30', '06:58', '14_ABASTO', 1, 'LP1912')
    55 = @('06:15:04', '07:16', '11_ETCHEVERRY', 61, 'LP1912')
    56 = @('06:43:40', '07:16', '16_SANTA ANA', 33, 'LP1912')
    117 = @('08:39:44', '09:34', '16_SANTA ANA', 55, 'LP1912')
    118 = @('08:39:44', '09:34', '23_HERNANDEZ', 55, 'LP1912')
    167 = @('11:20:07', '11:21', '16_SANTA ANA', 1, 'LP1912')
    168 = @('09:38:09', '11:21', '26_HERNANDEZ', 103, 'LP1912')
    169 = @('11:20:07', '11:22', '17_ROMERO', 2, 'LP1912')
    170 = @('10:56:30', '11:24', '10_OLMOS', 28, 'LP1912')
    171 = @('11:20:07', '11:25', '16_SANTA ANA', 5, 'LP1912')
    172 = @('09:38:09', '11:27', '225_C ROCA-H SUR', 109, 'LP1912')
    173 = @('09:38:09', '11:32', '81_EL PELIGRO', 114, 'LP1912')
    174 = @('10:56:30', '11:34', '23_HERNANDEZ', 38, 'LP1912')
    175 = @('10:26:41', '11:35', '11_ETCHEVERRY', 69, 'LP1912')
    176 = @('11:20:07', '11:35', '23_HERNANDEZ', 15, 'LP1912')
    177 = @('09:38:09', '11:36', '11_ETCHEVERRY', 118, 'LP1912')
    178 = @('10:26:41', '11:41', '17_ROMERO', 75, 'LP1912')
    179 = @('10:56:30', '11:42', '17_ROMERO', 46, 'LP1912')
    180 = @('10:26:41', '11:51', '215B_EL PATO', 85, 'LP1912')
    181 = @('10:56:30', '11:52', '15_ABASTO', 56, 'LP1912')
    182 = @('10:26:41', '11:59', '225_GOMEZ', 93, 'LP1912')
    183 = @('10:26:41', '12:02', '84_COLONIA URQUIZA-ESC 49', 96, 'LP1912')
    184 = @('11:20:07', '12:05', '23_HERNANDEZ', 45, 'LP1912')
    185 = @('10:26:41', '12:06', '16_P MOR-SANTA ANA', 100, 'LP1912')
    186 = @('10:56:30', '12:06', '14_ABASTO', 70, 'LP1912')
    187 = @('11:20:07', '12:07', '14_ABASTO', 47, 'LP1912')
    188 = @('11:20:07', '12:07', '16_P MOR-SANTA ANA', 47, 'LP1912')
    189 = @('10:56:30', '12:10', '10_OLMOS', 74, 'LP1912')
    190 = @('11:20:07', '12:13', '10_OLMOS', 53, 'LP1912')
    191 = @('10:26:41', '12:14', '17_ROMERO', 108, 'LP1912')
    192 = @('10:26:41', '12:19', '14_ABASTO', 113, 'LP1912')
    193 = @('10:26:41', '12:20', '215A_EL PATO', 114, 'LP1912')
    194 = @('10:56:30', '12:20', '14_ABASTO', 84, 'LP1912')
    195 = @('10:26:41', '12:21', '26_HERNANDEZ', 115, 'LP1912')
    196 = @('11:20:07', '12:21', '14_ABASTO', 61, 'LP1912')
    197 = @('11:20:07', '12:21', '215A_EL PATO', 61, 'LP1912')
    198 = @('10:56:30', '12:36', '27_EL RETIRO', 100, 'LP1912')
    199 = @('11:20:07', '12:37', '27_EL RETIRO', 77, 'LP1912')
    200 = @('10:56:30', '12:38', '17_179 Y 38', 102, 'LP1912')
    201 = @('10:56:30', '12:41', '10_OLMOS', 105, 'LP1912')
    202 = @('11:20:07', '12:49', '11_ETCHEVERRY', 89, 'LP1912')
    203 = @('11:20:07', '13:02', '15_ABASTO', 102, 'LP1912')
    204 = @('11:20:07', '13:07', '16_P MOR-SANTA ANA', 107, 'LP1912')
    205 = @('11:20:07', '13:14', '215D_EL PATO', 114, 'LP1912')
}
foreach ($r in $ws1Rows.Keys) {
    $vals = $ws1Rows[$r]
    $rn = [int]$r
    $ws1.Cells.Item($rn, 1).Value = $vals[0]
    $ws1.Cells.Item($rn, 2).Value = $vals[1]
    $ws1.Cells.Item($rn, 3).Value = $vals[2]
    $ws1.Cells.Item($rn, 4).Value = $vals[3]
    $ws1.Cells.Item($rn, 5).Value = $vals[4]
}

# ---- Sheet: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")

# Header updates
$ws2.Cells.Item(2, 1).Value = "Última actualización: 11:20:08"
$ws2.Cells.Item(3, 1).Value = "Total filas: 25"

# New rows appended
$ws2Rows = @{
    29 = @('11:20:07', '12:21', '215A_EL PATO', 61, 'LP1912')
    30 = @('11:20:07', '13:14', '215D_EL PATO', 114, 'LP1912')
}
foreach ($r in $ws2Rows.Keys) {
    $vals = $ws2Rows[$r]
    $rn = [int]$r
    $ws2.Cells.Item($rn, 1).Value = $vals[0]
    $ws2.Cells.Item($rn, 2).Value = $vals[1]
    $ws2.Cells.Item($rn, 3).Value = $vals[2]
    $ws2.Cells.Item($rn, 4).Value = $vals[3]
    $ws2.Cells.Item($rn, 5).Value = $vals[4]
}

# ---- Sheet: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2, 1).Value = "Última actualización: 11:20:08"

